# plantilla_localizaciones.xlsx — rename the single worksheet
# from "Hoja1" to "Ubicaciones" (per commit's accompanying diff).
#
# Note: the diff also shows Excel-build/theme-version metadata churn
# (fileVersion/rupBuild, absPath, revisionPtr GUID, window geometry,
# the default "Office" theme being renamed to "Office 2013 - 2022",
# and the consequent bestFit column-width / default-row-height drift).
# Those are side effects of which Excel build happened to resave the
# workbook, not an edit a user/script performs, so they are left to
# the host application to manage and are not reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Ubicaciones"
